# Fix typo in shared string (missing space), correct the "general"
# alignment of column B's blank cells to "left", widen column B so its
# content fits, and bump the row height of every data row from 18 to
# 18.75 pts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Resolve the data error: "buitenlandsepartijen" -> "buitenlandse partijen"
$ws.Range("B4").Value = "Ja, Nederlandse partijen; Nee, buitenlandse partijen"

# 2. The empty cells in column B were left-aligned as part of the fix
#    (they previously inherited the default "general" alignment).
$ws.Range("B1").HorizontalAlignment = -4131   # xlLeft
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B7").HorizontalAlignment = -4131
$ws.Range("B9").HorizontalAlignment = -4131

# 3. Column B is widened considerably so the (now longer) answers fit.
#    (74.3 is the ColumnWidth that Excel stores closest to the target
#    stored width of ~75.148 once it is rounded to whole pixels.)
$ws.Columns.Item(2).ColumnWidth = 74.3

# 4. Every data row (1-21) grows slightly taller.
for ($r = 1; $r -le 21; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}
